$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New column D values for the "practice" rows (2-5) ---
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# --- New column J values for the "generic" rows (6-9) ---
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# --- Fill in the new "unique_video" / "unique_audio" block, rows 14-21 ---
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "look"

$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "look"

$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "where"

$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "where"

$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "can"

$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "can"

$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "do"

$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "do"
